function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.281.47'
$ws.Range("E2").Value = '  +0.70%  '
$ws.Range("D3").Value = '1.831.59'
$ws.Range("E3").Value = '  -0.15%  '
Set-TextValue $ws.Range("D4") '0.9988'
$ws.Range("E4").Value = '  -0.19%  '
Set-TextValue $ws.Range("D5") '243.17'
$ws.Range("E5").Value = '  +0.43%  '
Set-TextValue $ws.Range("D6") '0.6183'
$ws.Range("E6").Value = '  -0.13%  '
Set-TextValue $ws.Range("D7") '0.9998'
$ws.Range("E7").Value = '  -0.16%  '
Set-TextValue $ws.Range("D8") '0.07364'
$ws.Range("E8").Value = '  -1.22%  '
Set-TextValue $ws.Range("D9") '0.2923'
$ws.Range("E9").Value = '  -0.11%  '
Set-TextValue $ws.Range("D10") '23.25'
$ws.Range("E10").Value = '  +0.61%  '
Set-TextValue $ws.Range("D11") '0.07649'
$ws.Range("E11").Value = '  -0.30%  '
$ws.Range("D12").Value = '1.854.23'
$ws.Range("E12").Value = '  +1.11%  '
Set-TextValue $ws.Range("D13") '4.994'
$ws.Range("E13").Value = '  -0.36%  '
Set-TextValue $ws.Range("D14") '0.6757'
$ws.Range("E14").Value = '  +0.20%  '
Set-TextValue $ws.Range("D15") '82.74'
$ws.Range("E15").Value = '  -0.37%  '
Set-TextValue $ws.Range("D16") '0.000008972'
$ws.Range("E16").Value = '  -1.93%  '
Set-TextValue $ws.Range("D17") '5.889'
$ws.Range("E17").Value = '  -0.46%  '
$ws.Range("D18").Value = '29.275.75'
$ws.Range("E18").Value = '  +0.74%  '
$ws.Range("D19").Value = '2.092.39'
$ws.Range("E19").Value = '  +0.40%  '
Set-TextValue $ws.Range("D20") '239.28'
$ws.Range("E20").Value = '  -0.82%  '
$ws.Range("E21").Value = '  -1.56%  '
Set-TextValue $ws.Range("D22") '0.9995'
$ws.Range("E22").Value = '  -0.19%  '
Set-TextValue $ws.Range("D23") '7.383'
$ws.Range("E23").Value = '  +2.36%  '
Set-TextValue $ws.Range("D24") '0.9999'
$ws.Range("E24").Value = '  -0.26%  '
Set-TextValue $ws.Range("D25") '158.43'
$ws.Range("E25").Value = '  -0.36%  '
$ws.Range("E26").Value = '  -1.06%  '
Set-TextValue $ws.Range("D27") '8.567'
Set-TextValue $ws.Range("D28") '17.68'
$ws.Range("E28").Value = '  -1.21%  '
Set-TextValue $ws.Range("D29") '1.494'
$ws.Range("E29").Value = '  -0.35%  '
Set-TextValue $ws.Range("D30") '0.05856'
$ws.Range("E30").Value = '  +4.13%  '
$ws.Range("E31").Value = '  +2.58%  '
Set-TextValue $ws.Range("D32") '4.084'
$ws.Range("E32").Value = '  -0.85%  '
Set-TextValue $ws.Range("D33") '4.099'
$ws.Range("E33").Value = '  -0.91%  '
Set-TextValue $ws.Range("D34") '1.860'
$ws.Range("E34").Value = '  +0.92%  '
Set-TextValue $ws.Range("D36") '0.7205'
$ws.Range("E36").Value = '  -2.88%  '
Set-TextValue $ws.Range("D37") '2.614'
$ws.Range("E37").Value = '  -1.64%  '
Set-TextValue $ws.Range("D38") '2.858'
$ws.Range("E38").Value = '  +3.05%  '
$ws.Range("D39").Value = '1.222.65'
$ws.Range("E39").Value = '  +0.96%  '
$ws.Range("E40").Value = '  -1.22%  '
Set-TextValue $ws.Range("D41") '0.9135'
$ws.Range("E41").Value = '  +1.52%  '
Set-TextValue $ws.Range("D42") '6.218'
$ws.Range("E42").Value = '  -2.90%  '
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("D44").Value = '2.010.58'
$ws.Range("E44").Value = '  +1.29%  '
Set-TextValue $ws.Range("D45") '101.88'
$ws.Range("E45").Value = '  +0.29%  '
Set-TextValue $ws.Range("D46") '65.83'
Set-TextValue $ws.Range("D47") '0.5056'
$ws.Range("E47").Value = '  -0.78%  '
$ws.Range("B48").Value = 'Algorand'
$ws.Range("C48").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range("D48") '0.1182'
$ws.Range("E48").Value = '  +6.69%  '
$ws.Range("B49").Value = 'TheSandbox'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range("D49") '0.4052'
$ws.Range("E49").Value = '  -0.31%  '
Set-TextValue $ws.Range("D50") '9.186'
$ws.Range("E50").Value = '  +0.72%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range("D51") '0.00000000117'
$ws.Range("E51").Value = '  -0.53%  '
